$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue 'D2' '30.463.78'
Set-TextValue 'E2' '  +0.27%  '

Set-TextValue 'D3' '1.866.43'
Set-TextValue 'E3' '  -0.31%  '

Set-TextValue 'E4' '  -0.05%  '

Set-TextValue 'D5' '235.36'
Set-TextValue 'E5' '  -0.91%  '

Set-TextValue 'E6' '  -0.02%  '

Set-TextValue 'D7' '0.4827'
Set-TextValue 'E7' '  +0.15%  '

Set-TextValue 'D8' '0.2809'
Set-TextValue 'E8' '  -0.27%  '

Set-TextValue 'D9' '0.06508'
Set-TextValue 'E9' '  -0.39%  '

Set-TextValue 'D10' '1.886.15'
Set-TextValue 'E10' '  +0.58%  '

Set-TextValue 'D11' '0.07437'
Set-TextValue 'E11' '  -0.10%  '

Set-TextValue 'D12' '16.40'
Set-TextValue 'E12' '  -0.28%  '

Set-TextValue 'D13' '5.061'
Set-TextValue 'E13' '  -0.34%  '

Set-TextValue 'D14' '87.31'
Set-TextValue 'E14' '  -0.72%  '

Set-TextValue 'D15' '0.6470'
Set-TextValue 'E15' '  -1.34%  '

Set-TextValue 'D16' '30.443.84'
Set-TextValue 'E16' '  +0.13%  '

Set-TextValue 'E17' '  +0.06%  '

$ws.Range('B18').Value = 'Avalanche'
$ws.Range('C18').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D18' '13.00'
Set-TextValue 'E18' '  -2.26%  '

$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D19' '234.01'
Set-TextValue 'E19' '  +5.79%  '

Set-TextValue 'D20' '0.000007532'
Set-TextValue 'E20' '  -1.20%  '

Set-TextValue 'D21' '2.110.12'
Set-TextValue 'E21' '  -0.35%  '

Set-TextValue 'D22' '1.001'
Set-TextValue 'E22' '  -0.08%  '

Set-TextValue 'D23' '5.151'
Set-TextValue 'E23' '  -2.62%  '

$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D24' '6.101'
Set-TextValue 'E24' '  -1.39%  '

$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D25' '9.342'
Set-TextValue 'E25' '  +1.09%  '

$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D26' '167.42'
Set-TextValue 'E26' '  +1.45%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D27' '18.40'
Set-TextValue 'E27' '  -0.68%  '

$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D28' '1.924'
Set-TextValue 'E28' '  -2.62%  '

$ws.Range('B29').Value = 'Stellar'
$ws.Range('C29').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D29' '0.1030'
Set-TextValue 'E29' '  +9.58%  '

$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D30' '1.374'
Set-TextValue 'E30' '  -5.57%  '

$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D31' '4.272'
Set-TextValue 'E31' '  -0.61%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D32' '4.017'
Set-TextValue 'E32' '  +0.06%  '

$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D33' '0.04984'
Set-TextValue 'E33' '  -1.03%  '

$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D34' '1.179'
Set-TextValue 'E34' '  -2.68%  '

$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D35' '0.7460'
Set-TextValue 'E35' '  -1.21%  '

$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 'D36' '1.000'
Set-TextValue 'E36' '  +0.15%  '

$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D37' '2.712'
Set-TextValue 'E37' '  +0.39%  '

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D38' '0.01934'
Set-TextValue 'E38' '  +5.23%  '

$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D39' '2.635'
Set-TextValue 'E39' '  +0.53%  '

$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D40' '0.9178'
Set-TextValue 'E40' '  +1.19%  '

$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D41' '2.055'
Set-TextValue 'E41' '  -1.21%  '

$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D42' '106.22'
Set-TextValue 'E42' '  -0.49%  '

$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D43' '0.9962'
Set-TextValue 'E43' '  -0.72%  '

$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D44' '0.4203'
Set-TextValue 'E44' '  -1.90%  '

$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D45' '5.541'
Set-TextValue 'E45' '  -6.73%  '

$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D46' '7.243'
Set-TextValue 'E46' '  -2.63%  '

$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D47' '61.91'
Set-TextValue 'E47' '  -5.56%  '

$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D48' '0.1231'
Set-TextValue 'E48' '  -5.34%  '

$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D49' '8.854'
Set-TextValue 'E49' '  -1.43%  '

$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D50' '1.441'
Set-TextValue 'E50' '  -2.30%  '

$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue 'D51' '33.62'
Set-TextValue 'E51' '  -1.57%  '
